$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '60.866.61'
$ws.Range("E2").Value = '  +0.92%  '

# Row 3
$ws.Range("D3").Value = '2.345.94'
$ws.Range("E3").Value = '  -0.85%  '

# Row 4
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''544.70'
$ws.Range("E5").Value = '  +1.01%  '

# Row 6
$ws.Range("D6").Value = '''136.70'
$ws.Range("E6").Value = '  -1.73%  '

# Row 7
$ws.Range("E7").Value = '  +0.04%  '

# Row 8
$ws.Range("D8").Value = '''0.525'
$ws.Range("E8").Value = '  -8.44%  '

# Row 9
$ws.Range("D9").Value = '2.340.02'
$ws.Range("E9").Value = '  -1.10%  '

# Row 10
$ws.Range("D10").Value = '''0.104'
$ws.Range("E10").Value = '  +0.30%  '

# Row 11
$ws.Range("E11").Value = '  +1.76%  '

# Row 12
$ws.Range("D12").Value = '''5.29'
$ws.Range("E12").Value = '  -0.17%  '

# Row 13
$ws.Range("D13").Value = '''0.341'
$ws.Range("E13").Value = '  +0.52%  '

# Row 14
$ws.Range("D14").Value = '''24.59'
$ws.Range("E14").Value = '  -2.17%  '

# Row 15
$ws.Range("D15").Value = '2.771.33'
$ws.Range("E15").Value = '  -0.77%  '

# Row 16
$ws.Range("D16").Value = '60.875.15'
$ws.Range("E16").Value = '  +1.24%  '

# Row 17
$ws.Range("D17").Value = '''0.0000159'
$ws.Range("E17").Value = '  -1.80%  '

# Row 18
$ws.Range("D18").Value = '2.349.08'
$ws.Range("E18").Value = '  -0.83%  '

# Row 19
$ws.Range("D19").Value = '''10.60'
$ws.Range("E19").Value = '  +0.76%  '

# Row 20
$ws.Range("D20").Value = '''318.89'
$ws.Range("E20").Value = '  +0.90%  '

# Row 21
$ws.Range("D21").Value = '''4.11'
$ws.Range("E21").Value = '  +1.10%  '

# Row 22
$ws.Range("D22").Value = '''6.52'
$ws.Range("E22").Value = '  -1.64%  '

# Row 23
$ws.Range("E23").Value = '  +0.10%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''63.33'
$ws.Range("E24").Value = '  +0.87%  '

# Row 25
$ws.Range("B25").Value = 'SuiNetwork'
$ws.Range("C25").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("D25").Value = '''1.72'
$ws.Range("E25").Value = '  -4.13%  '

# Row 26
$ws.Range("D26").Value = '''8.48'
$ws.Range("E26").Value = '  +10.86%  '

# Row 27
$ws.Range("D27").Value = '2.461.67'
$ws.Range("E27").Value = '  -0.94%  '

# Row 28
$ws.Range("D28").Value = '''7.94'
$ws.Range("E28").Value = '  +0.10%  '

# Row 29
$ws.Range("D29").Value = '''496.38'
$ws.Range("E29").Value = '  -4.05%  '

# Row 30
$ws.Range("D30").Value = '''1.37'
$ws.Range("E30").Value = '  -2.27%  '

# Row 31
$ws.Range("B31").Value = 'Kaspa'
$ws.Range("C31").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D31").Value = '''0.146'
$ws.Range("E31").Value = '  +2.14%  '

# Row 32
$ws.Range("B32").Value = 'PEPE'
$ws.Range("C32").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D32").Value = '0.0₃0859'
$ws.Range("E32").Value = '  -6.55%  '

# Row 33
$ws.Range("E33").Value = '  -1.50%  '

# Row 34
$ws.Range("D34").Value = '''1.49'
$ws.Range("E34").Value = '  -3.81%  '

# Row 35
$ws.Range("D35").Value = '''0.999'
$ws.Range("E35").Value = '  +0.14%  '

# Row 36
$ws.Range("D36").Value = '''4.61'
$ws.Range("E36").Value = '  +0.32%  '

# Row 37
$ws.Range("D37").Value = '''0.376'
$ws.Range("E37").Value = '  +1.17%  '

# Row 38
$ws.Range("D38").Value = '''18.50'
$ws.Range("E38").Value = '  +3.05%  '

# Row 39
$ws.Range("B39").Value = 'RenderToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range("D39").Value = '''5.25'
$ws.Range("E39").Value = '  -2.88%  '

# Row 40
$ws.Range("B40").Value = 'Stacks'
$ws.Range("C40").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D40").Value = '''1.82'
$ws.Range("E40").Value = '  +7.05%  '

# Row 41
$ws.Range("D41").Value = '''142.03'
$ws.Range("E41").Value = '  +3.73%  '

# Row 42
$ws.Range("E42").Value = '  -0.07%  '

# Row 43
$ws.Range("D43").Value = '''40.53'
$ws.Range("E43").Value = '  +1.12%  '

# Row 44
$ws.Range("D44").Value = '''142.37'
$ws.Range("E44").Value = '  +2.50%  '

# Row 45
$ws.Range("D45").Value = '''3.55'
$ws.Range("E45").Value = '  +1.44%  '

# Row 46
$ws.Range("D46").Value = '''2.04'
$ws.Range("E46").Value = '  -5.93%  '

# Row 47
$ws.Range("D47").Value = '''0.0516'
$ws.Range("E47").Value = '  +0.83%  '

# Row 48
$ws.Range("D48").Value = '''19.02'
$ws.Range("E48").Value = '  -4.88%  '

# Row 49
$ws.Range("E49").Value = '  -0.63%  '

# Row 50
$ws.Range("D50").Value = '''0.0901'
$ws.Range("E50").Value = '  -2.04%  '

# Row 51
$ws.Range("E51").Value = '  -0.70%  '
